$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting the old row 3 (MuSCs target) down to row 4.
$ws.Rows.Item(3).Insert()

# --- Row 2: Target cluster becomes "ECs" (a new cluster introduced by the
# updated TPM data); the recalculated expression/specificity numbers change
# too.
$ws.Range("D2").Value = "ECs"
$ws.Range("M2").Value = 1.536685
$ws.Range("N2").Value = 4.610055
$ws.Range("O2").Value = 0.4822880013826122
$ws.Range("P2").Value = 0.4822880013826122
$ws.Range("Q2").Value = 0.18527196371
$ws.Range("R2").Value = 1.66744767339
$ws.Range("S2").Value = 0.4822880013826122
$ws.Range("T2").Value = 0.4822880013826122

# --- Row 3 (newly inserted, blank): Sending=MuSCs Ligand=Ostn Receptor=Npr3
# Target=FAPs, carrying what used to be row 2's numbers (pre-update) plus
# freshly recalculated specificity values.
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Ostn"
$ws.Range("C3").Value = "Npr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.120566
$ws.Range("H3").Value = 0.361698
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.195417
$ws.Range("N3").Value = 3.586251
$ws.Range("O3").Value = 0.3751811696924212
$ws.Range("P3").Value = 0.3751811696924212
$ws.Range("Q3").Value = 0.144126646022
$ws.Range("R3").Value = 1.297139814198
$ws.Range("S3").Value = 0.3751811696924212
$ws.Range("T3").Value = 0.3751811696924212

# --- Row 4 (formerly row 3, shifted down by the insert): Target stays
# "MuSCs"; only the derived-specificity columns (O/P/S/T) are recalculated
# because of the new cluster row, the rest carries over unchanged.
$ws.Range("O4").Value = 0.1425308289249667
$ws.Range("P4").Value = 0.1425308289249667
$ws.Range("S4").Value = 0.1425308289249667
$ws.Range("T4").Value = 0.1425308289249667
